# config.xlsx - "Testing, updated form, added breakpoint 5, fixed sep certs logic"
#
# 1. Insert a new "AddressLabelPath" parameter row just above "CountryLookUpPath"
#    (old row 30 -> new row 31), pushing every following row down by one.
# 2. Grow Table1 (the parameters table on Sheet1) by the one extra row.
# 3. Fix the "SeparateProducts" regex: the old pattern moves to the Description
#    column and a corrected pattern takes its place in the Value column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# --- 1. New row: AddressLabelPath -------------------------------------------
$ws.Rows.Item(30).Insert()

$ws.Cells.Item(30, 1).Value = "AddressLabelPath"
$ws.Cells.Item(30, 2).Value = "\\EARTH.GSI.GOV.UK\USER\SHARED\Agency\CoFS for G drive\RobotDocuments\RobotAddressLabelTemplate.docx"
$ws.Cells.Item(30, 3).Clear()
$ws.Rows.Item(30).RowHeight = 48.75

# --- 2. Table1 now spans one more row ---------------------------------------
$lo.Resize($ws.Range("A1:C56"))

# --- 3. Fixed regex for SeparateProducts (now row 53) -----------------------
$ws.Cells.Item(53, 3).Value = "{\Wproduct\W:\W(\d+\s*)+\W}"
$ws.Cells.Item(53, 2).Value = "{\Wproduct\W:\W((\w+\s*\W)+)}"

# --- View state (best effort) ------------------------------------------------
$ws.Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B57").Select()
